$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.273.72"
$ws.Range("E2").Value = "  +4.14%  "
$ws.Range("D3").Value = "2.430.05"
$ws.Range("E3").Value = "  +3.25%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "554.83"
$ws.Range("E5").Value = "  +2.19%  "
$ws.Range("E6").Value = "  +3.67%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").Value = "0.577"
$ws.Range("E8").Value = "  +1.31%  "
$ws.Range("E9").Value = "  +4.54%  "
$ws.Range("D10").Value = "5.76"
$ws.Range("E10").Value = "  +4.14%  "
$ws.Range("E11").Value = "  +0.64%  "
$ws.Range("E12").Value = "  -2.10%  "
$ws.Range("D13").Value = "25.00"
$ws.Range("E13").Value = "  +5.44%  "
$ws.Range("D14").Value = "2.863.11"
$ws.Range("E14").Value = "  +3.24%  "
$ws.Range("D15").Value = "60.204.12"
$ws.Range("E15").Value = "  +4.12%  "
$ws.Range("E16").Value = "  +4.15%  "
$ws.Range("D17").Value = "2.434.06"
$ws.Range("E17").Value = "  +3.14%  "
$ws.Range("E18").Value = "  +6.14%  "
$ws.Range("E19").Value = "  +3.17%  "
$ws.Range("D20").Value = "332.94"
$ws.Range("E20").Value = "  +0.82%  "
$ws.Range("D21").Value = "6.78"
$ws.Range("E21").Value = "  +1.07%  "
$ws.Range("E22").Value = "  +0.03%  "
$ws.Range("D23").Value = "65.19"
$ws.Range("E23").Value = "  +4.30%  "
$ws.Range("E24").Value = "  +3.43%  "
$ws.Range("D25").Value = "8.64"
$ws.Range("E25").Value = "  +2.80%  "
$ws.Range("E26").Value = "  -0.46%  "
$ws.Range("D27").Value = "1.35"
$ws.Range("E27").Value = "  -0.62%  "
$ws.Range("D28").Value = "0.0₃0788"
$ws.Range("E28").Value = "  +7.18%  "
$ws.Range("E29").Value = "  +1.47%  "
$ws.Range("D30").Value = "6.34"
$ws.Range("E30").Value = "  +3.41%  "
$ws.Range("D31").Value = "169.45"
$ws.Range("E31").Value = "  -0.60%  "
$ws.Range("E32").Value = "  +3.36%  "
$ws.Range("D33").Value = "18.74"
$ws.Range("E33").Value = "  +1.95%  "
$ws.Range("E34").Value = "  +0.01%  "
$ws.Range("E35").Value = "  +6.04%  "
$ws.Range("D36").Value = "4.23"
$ws.Range("E36").Value = "  -0.04%  "
$ws.Range("E37").Value = "  +0.09%  "
$ws.Range("E38").Value = "  +0.35%  "
$ws.Range("D39").Value = "325.17"
$ws.Range("E39").Value = "  +12.55%  "
$ws.Range("D40").Value = "0.421"
$ws.Range("E40").Value = "  +11.46%  "
$ws.Range("D41").Value = "39.76"
$ws.Range("E41").Value = "  +1.81%  "
$ws.Range("E42").Value = "  +1.59%  "
$ws.Range("D43").Value = "140.56"
$ws.Range("E43").Value = "  -1.46%  "
$ws.Range("E44").Value = "  +3.75%  "
$ws.Range("D45").Value = "0.0961"
$ws.Range("E45").Value = "  +1.18%  "
$ws.Range("D46").Value = "19.60"
$ws.Range("E46").Value = "  +2.17%  "
$ws.Range("E47").Value = "  +8.90%  "
$ws.Range("D48").Value = "0.574"
$ws.Range("E48").Value = "  +1.33%  "
$ws.Range("E49").Value = "  +2.13%  "
$ws.Range("D50").Value = "17.91"
$ws.Range("E50").Value = "  +2.53%  "
$ws.Range("D51").Value = "11.06"
$ws.Range("E51").Value = "  -0.34%  "
